$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural changes --------------------------------------------------
# Insert a new column before the old "Email" column (C) so the layout
# becomes Name | Date | Month | Email.
$ws.Columns("C").Insert()

# Insert a new header row at the top; existing rows shift down by one.
$ws.Rows("1").Insert()

# --- Header row ------------------------------------------------------------
$ws.Range("A1").Value = "Name"
$ws.Range("D1").Value = "Email"

# --- Data rows (Name / Email) -------------------------------------------
$ws.Range("A2").Value = "Avik Deb"
$ws.Range("D2").Value = "avikdeb@gmail.com"

$ws.Range("A3").Value = "Ayutasouri Deb"
$ws.Range("D3").Value = "ayutasouri@gmail.com"

$ws.Range("A4").Value = "Rabindranath Tagore"
$ws.Range("D4").Value = "avik.consult@gmail.com"

$ws.Range("A5").Value = "Mahatma Gandhi"
$ws.Range("D5").Value = "avik.consult@gmail.com"

$ws.Range("A6").Value = "Test Human"
$ws.Range("D6").Value = "avikdeb@gmail.com"

# --- Month column (C): always text -----------------------------------
$ws.Range("C1:C6").NumberFormat = "@"
$ws.Range("C1").Value = "Month"
$ws.Range("C2").Value = "Apr"
$ws.Range("C3").Value = "Sep"
$ws.Range("C4").Value = "May"
$ws.Range("C5").Value = "Oct"
$ws.Range("C6").Value = "Apr"

# --- Date column (B) ------------------------------------------------------
# The column is Text-formatted, but only the header + first data row were
# actually typed as text; the remaining day-of-month numbers are stored
# as real numbers (format applied after the value is entered keeps the
# numeric storage while still displaying/validating as text).
$ws.Range("B1").NumberFormat = "@"
$ws.Range("B1").Value = "Date"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "1"

$ws.Range("B3").Value = 23
$ws.Range("B3").NumberFormat = "@"

$ws.Range("B4").Value = 7
$ws.Range("B4").NumberFormat = "@"

$ws.Range("B5").Value = 2
$ws.Range("B5").NumberFormat = "@"

$ws.Range("B6").Value = 10
$ws.Range("B6").NumberFormat = "@"

# --- Column widths -------------------------------------------------------
$ws.Columns("B").ColumnWidth = 4.3
$ws.Columns("C").ColumnWidth = 6
$ws.Columns("D").ColumnWidth = 22.2

# --- Sheet view / page setup ----------------------------------------------
$ws.Range("D6").Select() | Out-Null
$ws.PageSetup.Orientation = 1
